$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = "[Null space control] redundant manipulator control with optimal target function"
$ws.Range("E28").Value = "https://ropiens.tistory.com/152"

$ws.Range("D33").Value = "Tabular Data(정형 데이터)에서의 Noise"
$ws.Range("E33").Value = "https://velog.io/@vvakki_/Tabular-Data%EC%A0%95%ED%98%95-%EB%8D%B0%EC%9D%B4%ED%84%B0%EC%97%90%EC%84%9C%EC%9D%98-Noise"

$ws.Range("D44").Value = "Meta-Learning Based Beamforming 논문 리뷰"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/96"
